# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 505

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = 45181
}
